{"js": "// \"Getting rid of AR\" \u2014 replace the \"AR\" shorthand with \"Recommendation\"\n// wherever it appears as the document's own terminology (title + body\n// sentence). The merge-field ${AR} placeholder becomes ${REC} to match.\n\n// 1) Title: \"AR ${AR}: Switch to LED Lighting\" -> \"Recommendation ${REC}: Switch to LED Lighting\"\nconst titleResults = context.document.body.search(\"AR ${AR}\", { matchCase: true, matchWholeWord: false });\ntitleResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < titleResults.items.length; i++) {\n  titleResults.items[i].insertText(\"Recommendation ${REC}\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Body sentence: \"This AR will indicate\" -> \"This recommendation will indicate\"\nconst bodyResults = context.document.body.search(\"This AR will indicate\", { matchCase: true, matchWholeWord: false });\nbodyResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < bodyResults.items.length; i++) {\n  bodyResults.items[i].insertText(\"This recommendation will indicate\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# \"Getting rid of AR\" \u2014 replace the \"AR\" shorthand with \"Recommendation\"\n# wherever it appears as the document's own terminology (title + body\n# sentence). The merge-field ${AR} placeholder becomes ${REC} to match.\n#\n# NOTE: use single-quoted strings for any literal containing \"${...}\" so\n# PowerShell does not try to expand it as a variable reference.\n\n$d = $word.ActiveDocument\n\n# 1) Title: \"AR ${AR}: Switch to LED Lighting\" -> \"Recommendation ${REC}: Switch to LED Lighting\"\n$find1 = $d.Content.Find\n$find1.Text = 'AR ${AR}'\n$find1.Replacement.Text = 'Recommendation ${REC}'\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2)\n\n# 2) Body sentence: \"This AR will indicate\" -> \"This recommendation will indicate\"\n$find2 = $d.Content.Find\n$find2.Text = 'This AR will indicate'\n$find2.Replacement.Text = 'This recommendation will indicate'\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n"}
